# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric to Excel auto-detection (e.g. "238.25"),
# but the source data is text (prices formatted with literal dot separators,
# e.g. "40.724.03" uses dots as thousands separators). Force text storage
# via a temporary Text number format, then clear the format back off so the
# cell keeps the default (unstyled) appearance, matching the original file.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "40.724.03"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -7.06%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.180.63"
$cell.ClearFormats()
$ws.Range("E3").Value = "  -7.47%  "

$ws.Range("E4").Value = "  -0.02%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "238.25"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.70%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.617"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -7.86%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "69.54"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -5.22%  "

$ws.Range("E8").Value = "  +0.16%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.532"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -11.44%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "36.29"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +3.81%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "57.48"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -5.58%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.0931"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -8.83%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -4.48%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.49"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -10.01%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.504.71"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -7.58%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "14.46"
$cell.ClearFormats()
$ws.Range("E16").Value = "  -10.51%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.825"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -9.30%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.203.78"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -6.72%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "40.728.96"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -7.05%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0931"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -9.71%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "71.83"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -7.53%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.00"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -8.25%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "229.43"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -9.12%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.99"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +6.88%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "3.58"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -5.02%  "

$ws.Range("E27").Value = "  -4.65%  "

$ws.Range("E28").Value = "  -5.23%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.62"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -8.17%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "168.51"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -4.24%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "20.04"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -9.94%  "

$ws.Range("E32").Value = "  -10.00%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.122"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -8.24%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0694"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -7.32%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.03"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -5.71%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "4.52"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -10.23%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "3.78"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -0.18%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "22.81"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +12.63%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -7.24%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.0263"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -4.62%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.77"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -12.30%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "63.83"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -0.95%  "

$ws.Range("E43").Value = "  -11.17%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "8.59"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -4.87%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.190"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -6.36%  "

$ws.Range("E46").Value = "  -0.14%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.0975"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -7.67%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "4.42"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +1.53%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "10.01"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +5.29%  "

$ws.Range("E50").Value = "  -6.44%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.07"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -6.91%  "
